$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: column A = cycles/temp values, column B = old A values, column C = new constant 60
$ws.Range("A1").Value = 5000
$ws.Range("B1").Value = 105
$ws.Range("C1").Value = 60

$ws.Range("A2").Value = 10000
$ws.Range("B2").Value = 105
$ws.Range("C2").Value = 60

$ws.Range("A3").Value = 12500
$ws.Range("B3").Value = 150
$ws.Range("C3").Value = 60

$ws.Range("A4").Value = 20000
$ws.Range("B4").Value = 150
$ws.Range("C4").Value = 60

$ws.Range("A5").Value = 25000
$ws.Range("B5").Value = 30
$ws.Range("C5").Value = 60

# Update the active selection to match the target (A6 instead of B6)
$ws.Range("A6").Select()
